$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Cell F27 already contains the literal text "TRUE" (as a shared string,
# not the Boolean value). Copy it and paste-special (values only) into
# each target cell so the destination keeps its own formatting/style
# while its content becomes the text "TRUE" rather than being
# auto-coerced into a Boolean by a plain .Value assignment.
$source = $ws.Range("F27")
[void]$source.Copy()

$targetAddresses = @(
    "E4", "E5", "E6",
    "E8", "E9", "E10", "E11", "E12", "E13",
    "E15",
    "E17",
    "E19", "E20", "E21", "E22",
    "E24", "E25", "E26",
    "E29", "E30", "E31", "E32", "E33"
)

foreach ($addr in $targetAddresses) {
    [void]$ws.Range($addr).PasteSpecial(-4163, -4142, $false, $false)
}

$excel.CutCopyMode = 0

# Update the selected/active cell as recorded in the sheet view.
[void]$ws.Range("E34").Select()
